$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A228").Value = "Jackley clarifies the focus of multi-state Section 504 lawsuit"
$ws.Range("A229").Value = "Lawsuit filed against Keller ISD alleges violations of Voting Rights Act"
$ws.Range("A230").Value = "20 red states, including Arkansas, back DOGE in lawsuit challenging access to Treasury system"
$ws.Range("A231").Value = "Prof. James Sample Explains Federal Lawsuit Against New York’s Green Light Law"
$ws.Range("A232").Value = "Law’s Sharona Hoffman discusses an insulin lawsuit recently filed by the City of Columbus"
$ws.Range("A233").Value = "Elmore County lawyer arrested for harassment"
$ws.Range("A234").Value = "Fla. Lawyer Allegedly Smashed Plate Over Fellow Wedding Attendee’s Head When He Allowed Others to Cut in Buffet Line"
$ws.Range("A235").Value = "Attorney general threatening legal action against Indianapolis authorities, schools"
$ws.Range("A236").Formula = "=""'Outrageous and false': How those named in the AHS lawsuit are responding to the allegations"""
$ws.Range("A236").Copy()
$ws.Range("A236").PasteSpecial(-4163)
$ws.Range("A237").Value = "GloRilla Denies BBL Rumors as She Preps for ‘The Glorious Tour’"
$ws.Range("A238").Value = "ED seizes Rs 170-cr worth bank deposits in probe against ‘fraud’ forex trading platform"
$ws.Range("A239").Value = "DOGE Posts — Then Redacts — What Appears to Be Sensitive HUD Contract Data"
$ws.Range("A240").Value = "Tarver Elementary teacher placed on leave amid misconduct allegations - KWKT"
$ws.Range("A241").Value = "Mobile Mardi Gras queen accused of nearly `$1.5M embezzlement scheme"
$ws.Range("A242").Value = "Vernon business owner accused of child sex crimes - KFDX"
$ws.Range("A243").Value = "Mass. daycare co-owner accused of assaulting children in her care"
$ws.Range("A244").Value = "B1/B2 Visa Revoked -Traveler Finds Out at Airport"
$ws.Range("A245").Value = "EHarley Street: Calls for inquiry into GP management 'scandal'"
$ws.Range("A246").Value = "Netflix Addresses Controversy Surrounding Karla Sofía Gascón Amid ‘Emilia Pérez’ Scandal"
$excel.CutCopyMode = 0
